$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($r in 9..11) {
    $ws.Cells.Item($r, 1).Value = "JERSH"
    $ws.Cells.Item($r, 2).Value = "MERXWERLS"
    $ws.Cells.Item($r, 4).Value = "123-111-9928"
}
